$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell as TEXT (matching the
# workbook's existing inline-string / text-typed Price column) without
# leaving a residual NumberFormat/style change on the cell. A leading
# apostrophe forces Excel to store the value as text; resetting the style
# back to "Normal" afterwards clears the quote-prefix style Excel applies
# so the cell's style index is left exactly as it was before.
function Set-TextValue($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "243.96"
Set-TextValue "D3"  "23.83"
Set-TextValue "D4"  "5.242"
Set-TextValue "D5"  "0.05814"
Set-TextValue "D6"  "6.466"
Set-TextValue "D7"  "3.229"
Set-TextValue "D8"  "0.8082"
Set-TextValue "D9"  "0.8820"
Set-TextValue "D10" "0.1393"
Set-TextValue "D11" "0.07096"
Set-TextValue "D12" "0.03285"
Set-TextValue "D14" "0.09329"
Set-TextValue "D15" "3.832"
Set-TextValue "D16" "0.001551"
Set-TextValue "D17" "0.04709"
Set-TextValue "D18" "0.0006012"
Set-TextValue "D19" "0.006156"
Set-TextValue "D20" "0.001259"
Set-TextValue "D21" "0.004071"
Set-TextValue "D22" "0.00008704"
Set-TextValue "D24" "2.153"
Set-TextValue "D25" "0.3185"
Set-TextValue "D26" "0.1321"
Set-TextValue "D28" "0.0002329"
Set-TextValue "D40" "0.03789"
Set-TextValue "D44" "0.007858"
Set-TextValue "D45" "0.00005315"
Set-TextValue "D47" "0.5351"
Set-TextValue "D48" "0.002597"

# --- E18 text update (appended "Worstin24h" marker) ---
$ws.Range("E18").Value = "17OneONEWorstin24h"

# --- Rows 41-43 rotate: KickToken -> row41, BKEXToken -> row42, CEJI -> row43 ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006283"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1051"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002601"
$ws.Range("E43").Value = "42CEJICEJI"
